# Insert a new weekly price record as row 266 on the single worksheet,
# pushing the previously existing rows 266-327 down to 267-328.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(266).EntireRow.Insert()

$ws.Cells.Item(266, 1).Value  = 8
$ws.Cells.Item(266, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(266, 3).Value  = 'Coquimbo'
$ws.Cells.Item(266, 4).Value  = 44711
$ws.Cells.Item(266, 5).Value  = 4
$ws.Cells.Item(266, 6).Value  = 100112032
$ws.Cells.Item(266, 7).Value  = 'Zapallo italiano'
$ws.Cells.Item(266, 8).Value  = 'Sin especificar'
$ws.Cells.Item(266, 9).Value  = 'Primera'
$ws.Cells.Item(266, 10).Value = 440
$ws.Cells.Item(266, 11).Value = 14000
$ws.Cells.Item(266, 12).Value = 15000
$ws.Cells.Item(266, 13).Value = 14500
$ws.Cells.Item(266, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(266, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(266, 16).Value = 242
$ws.Cells.Item(266, 17).Value = 60
$ws.Cells.Item(266, 18).Value = 'Hortaliza'
